$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.854.75'
$ws.Range('E2').Value = '  -2.80%  '
$ws.Range('D3').Value = '2.566.22'
$ws.Range('E3').Value = '  -3.03%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '515.07'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.37'
$ws.Range('E6').Value = '  -5.02%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.558'
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('D9').Value = '2.580.87'
$ws.Range('E9').Value = '  -3.33%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.45'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0988'
$ws.Range('E11').Value = '  -5.41%  '
$ws.Range('E12').Value = '  -4.04%  '
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').Value = '3.023.98'
$ws.Range('E14').Value = '  -2.75%  '
$ws.Range('D15').Value = '57.859.06'
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.02'
$ws.Range('E16').Value = '  -4.20%  '
$ws.Range('D17').Value = '2.589.09'
$ws.Range('E17').Value = '  -5.09%  '
$ws.Range('E18').Value = '  -4.64%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '332.20'
$ws.Range('E19').Value = '  -3.30%  '
$ws.Range('E20').Value = '  -4.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.03'
$ws.Range('E21').Value = '  -5.88%  '
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.65'
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.164'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -4.89%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.89'
$ws.Range('E28').Value = '  -4.88%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0707'
$ws.Range('E30').Value = '  -12.05%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.89'
$ws.Range('E31').Value = '  -8.46%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.55'
$ws.Range('E32').Value = '  -4.18%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.57'
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '149.26'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.87'
$ws.Range('E35').Value = '  -7.47%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.11'
$ws.Range('E36').Value = '  -7.81%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '36.11'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.826'
$ws.Range('E38').Value = '  -5.55%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.825'
$ws.Range('E39').Value = '  -4.78%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.41'
$ws.Range('E40').Value = '  -5.62%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.49'
$ws.Range('E41').Value = '  -4.71%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '272.63'
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '10.69'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.589'
$ws.Range('E45').Value = '  -2.68%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0938'
$ws.Range('E46').Value = '  -4.04%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0513'
$ws.Range('E47').Value = '  -4.64%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.974.37'
$ws.Range('E48').Value = '  -3.21%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.33'
$ws.Range('E49').Value = '  -5.89%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.48'
$ws.Range('E50').Value = '  -6.23%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0217'
$ws.Range('E51').Value = '  -5.62%  '
